$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

# Update the username value in D2
$ws.Range("D2").Value = "extendhealth\piroy"

# Update the active selection to D2
$ws.Activate()
$ws.Range("D2").Select()
